$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (14) of data - data processing examples (groundwater model)
$ws.Range("A14").Value = "groundwater_model"
$ws.Range("B14").Value = "practical examples"
$ws.Range("D14").Value = "showcase"
$ws.Range("E14").Value = "flopy,  animation"
$ws.Range("F14").Value = "Onno Ebbens"
$ws.Range("H14").Value = "af"

# Match the wrap-text style used by the analogous cells in row 13
$ws.Range("B14").WrapText = $true
$ws.Range("E14").WrapText = $true
$ws.Range("H14").WrapText = $true
